$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Objetivos:): replace B/C with new objectives text ---
$objText = @'
O principal objetivo é permitir ao aluno conhecer os princípios fundamentais das interfaces líquido/gás/sólido e aprender sobre os conceitos de adsorção em sólidos, a caracterização de superfícies porosas, tensão superficial e propriedades de sistemas coloidais e emulsões. Além disso, identificar e explorar as aplicações destes conceitos em diferentes processos da indústria química.
'@
$ws.Range("B10").Value = $objText
$ws.Range("C10").Value = $objText

# --- Insert a new row at position 13; shifts old rows 13-24 down to 14-25 ---
$ws.Rows.Item(13).Insert()

# New row 13: only B13/C13 populated ("1488970 - Marivone Nunho Sousa"), no A13
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$marivoneText = @'
1488970 - Marivone Nunho Sousa
'@
$ws.Range("B13").Value = $marivoneText
$ws.Range("C13").Value = $marivoneText

# --- Row 14 (Programa resumido:): replace B/C with new description text ---
$descText = @'
Descrição de superfície e interface, termodinâmica das superfícies. Superfícies e forças. Sólidos iônicos e covalentes. Forças físicas e químicas de adsorção. Interface gás-solido e líquido-sólido. Caracterização de superfícies.
'@
$ws.Range("B14").Value = $descText
$ws.Range("C14").Value = $descText

# --- Row 16 (Programa:): replace B/C with new syllabus text ---
$programaText = @'
1)Isotermas de adsorção: Isotermas de Langmuir. Isotermas de Brunauer, Emmett e Teller. Métodos de determinação da área superficial de sólidos. Classificação quanto à porosidade. Métodos de determinação da porosidade de sólidos.2)Catálise de superfície. Interface sólido‐líquido. Mecanismos de catálise.3)Tensão superficial e interfacial. Equação de Laplace. Ângulo de contato. Ascensão e depressão capilar. Aplicações.4)Classificação das dispersões coloidais. Dupla camada elétrica: equação de Lippman e apresentação de modelos. Estabilidade e coagulação de dispersões coloidais. 5)Interações intermoleculares, dipolo-dipolo e de Van-der-Waals, ligação de hidrogênio e interações estabilizadoras em macromoléculas. 6)Estado coloidal. Colóides liofílicos e liofóbicos, hdrofílicos e hidrofóbicos. Obtenção de colóides. Propriedades cinéticas difusão, sedimentação, convecção. Propriedades óticas: espalhamento estático de luz, turbidez, espalhamento dinâmico da luz.Coagulação. Aplicações.7)Termodinâmica dos processos de transporte: difusão sedimentação e transporte através de membranas.
'@
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# --- Row 19 (Metodo:): B/C becomes the "Participacao..." text ---
$participacaoText = @'
Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas.
'@
$ws.Range("B19").Value = $participacaoText
$ws.Range("C19").Value = $participacaoText

# --- Row 20 (Criterio:): B/C becomes the "Media Final..." text ---
$mediaFinalText = @'
Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3
Média final mínima de aprovação = 5,0
'@
$ws.Range("B20").Value = $mediaFinalText
$ws.Range("C20").Value = $mediaFinalText

# --- Row 21 (Norma de recuperacao:): B/C becomes the "Prova escrita..." text ---
$provaEscritaText = @'
(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0
'@
$ws.Range("B21").Value = $provaEscritaText
$ws.Range("C21").Value = $provaEscritaText

# --- Row 22 (Bibliografia:): replace B/C with new bibliography text ---
$bibText = @'
1)MYERS, D. Surfaces, interfaces, and colloids: Principles and Applications, Second edition, Wiley-VCH, New York, 19982) BIRDI, K. S.; Surface and Colloid Chemistry, 1a ed., CRC Press LLC, New York, 1997.3) OSHIMA, H., Theory of colloid and interfacial electric phenomena. Interface Science and Technology Series, v. 12, Academic Press, Oxford, 2006.4) JACOB N. ISRAELCHVILI; Intermolecular and Surface Forces, 3r d Edition, New York, Academic, 2010.5) ADAMIAN, R. E ALMENDRA E.; Físico-Química – Uma Aplicação aos Materiais, 2002. 6) ADAMSON, A. Physical Chemistry of Surfaces (5th ed.). New York: John Wiley, 1990.7) SHAW, D. J. Introdução à Química dos Coloides e de Superfícies. São Paulo: Edgard Blücher, 1975. 185 pp.8)  REGALBUTO, J. Handbook of catalyst preparation. Taylor & Francis,2007
'@
$ws.Range("B22").Value = $bibText
$ws.Range("C22").Value = $bibText
